# edit.ps1 -- apply "Added NodeJS to CV." changes to the Technical Skills
# list, plus the accompanying bold-formatting fix on the empty paragraph
# right after the "EDUCATION" heading.
#
# Summary of the target change (see unified diff):
#   1. Add a new "Java" bullet right after "Ruby/JRuby" (before "F#").
#   2. Change the "CoffeeScript" bullet's text to "JavaScript/CoffeeScript"
#      (as two runs: "JavaScript/" + "CoffeeScript").
#   3. Remove the old "Java" bullet that used to sit between
#      "CoffeeScript" and "Groovy".
#   4. Add a new "NodeJS" bullet right before "Linux" (after "Groovy").
#   5. Make the empty paragraph directly under "EDUCATION" fully bold
#      (paragraph mark + run both get rFonts/b/sz/szCs).

function Find-ParaIndex {
    param($doc, [string]$text, [int]$startAt = 1)
    for ($i = $startAt; $i -le $doc.Paragraphs.Count; $i++) {
        $t = $doc.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert "Java" as a new bullet before "F#" (i.e. right after
#    "Ruby/JRuby"). InsertParagraphBefore() clones the paragraph mark
#    formatting (numbering, indents, fonts) from "F#", so the new
#    paragraph already has the correct list-item styling.
# ---------------------------------------------------------------------
$fSharpIdx = Find-ParaIndex $d "F#"
$d.Paragraphs($fSharpIdx).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs($fSharpIdx).Range.Text = "Java"

# ---------------------------------------------------------------------
# 2. Turn "CoffeeScript" into "JavaScript/" + "CoffeeScript" (two runs
#    with identical run formatting). Insert the new text at the start
#    of the paragraph, then force a run split by toggling Bold on the
#    inserted span (off again), which keeps the two pieces as separate
#    <w:r> runs instead of Word silently re-merging them.
# ---------------------------------------------------------------------
$coffeeIdx = Find-ParaIndex $d "CoffeeScript"
$coffeePara = $d.Paragraphs($coffeeIdx)
$insertionStart = $coffeePara.Range.Start
$collapsed = $d.Range($insertionStart, $insertionStart)
$collapsed.InsertBefore("JavaScript/")
$newRunRange = $d.Range($insertionStart, $insertionStart + 11)
$newRunRange.Bold = 1
$newRunRange.Bold = 0

# ---------------------------------------------------------------------
# 3. Remove the old "Java" bullet (now sitting between "CoffeeScript"
#    and "Groovy"). Search *after* the CoffeeScript paragraph so we
#    don't re-match the brand-new "Java" bullet added in step 1.
# ---------------------------------------------------------------------
$oldJavaIdx = Find-ParaIndex $d "Java" ($coffeeIdx + 1)
$d.Paragraphs($oldJavaIdx).Range.Delete() | Out-Null

# ---------------------------------------------------------------------
# 4. Insert "NodeJS" as a new bullet before "Linux" (i.e. right after
#    "Groovy").
# ---------------------------------------------------------------------
$linuxIdx = Find-ParaIndex $d "Linux"
$d.Paragraphs($linuxIdx).Range.InsertParagraphBefore() | Out-Null
$d.Paragraphs($linuxIdx).Range.Text = "NodeJS"

# ---------------------------------------------------------------------
# 5. Bold the empty paragraph directly below "EDUCATION". Setting Font
#    properties on the (collapsed) paragraph-mark range updates both
#    the run's rPr and the paragraph mark's rPr (w:pPr/w:rPr).
# ---------------------------------------------------------------------
$eduIdx = Find-ParaIndex $d "EDUCATION"
$afterEduRange = $d.Paragraphs($eduIdx + 1).Range
$afterEduRange.Font.Bold = 1
$afterEduRange.Font.Name = "Calibri"
$afterEduRange.Font.Size = 12
$afterEduRange.Font.SizeBi = 12

Write-Output "done"
